$wb = $excel.ActiveWorkbook
$ws5 = $wb.Worksheets.Item("Segment_Margins")
$rng = $ws5.Range("F2:F5")
$cs = $rng.FormatConditions.AddColorScale(2)
$cs.ColorScaleCriteria.Item(1).FormatColor.Color = 2181337
$wb.Save()
